$wb = $excel.ActiveWorkbook

# --- Cover sheet: add Table 5 / Table 6 caption rows ---
$wsCover = $wb.Worksheets.Item("Cover")
$wsCover.Cells.Item(6,1).Value = 'Table 5'
$wsCover.Cells.Item(6,2).Value = 'Differences in study variable distribution between participant with and without testis cancer relapse. Numeric variables are presented as medians with interqurtile ranges (IQR) and ranges. Categorical variables are presented as percentages and counts within the complete observation set.'
$wsCover.Cells.Item(7,1).Value = 'Table 6'
$wsCover.Cells.Item(7,2).Value = 'Characteristic of multi-paramater LASSO Cox models of effects of PRL concentration and PRL strata on relapse-free survival.'

# --- Table 1 sheet: append LDH_class / AFP_class / HCG_class rows ---
$wsT1 = $wb.Worksheets.Item("Table 1")
$wsT1.Cells.Item(46,1).Value = 'LDH_class'
$wsT1.Cells.Item(46,2).Value = 'Lactate dehydrogenase, blood concentration strata'
$wsT1.Cells.Item(46,3).Value = 'LDH strata'
$wsT1.Cells.Item(46,4).Value = '0 - 190 U/L, > 190 U/L'
$wsT1.Cells.Item(47,1).Value = 'AFP_class'
$wsT1.Cells.Item(47,2).Value = 'Alpha fetoprotein, blood concentration strata'
$wsT1.Cells.Item(47,3).Value = 'AFP strata'
$wsT1.Cells.Item(47,4).Value = '0 - 40 ng/mL, > 40 ng/mL'
$wsT1.Cells.Item(48,1).Value = 'HCG_class'
$wsT1.Cells.Item(48,2).Value = 'Human chorionic gonadotropin, blood concentration strata'
$wsT1.Cells.Item(48,3).Value = 'HCG strata'
$wsT1.Cells.Item(48,4).Value = '0 - 5 IU/L, > 5 IU/L'

# --- Table 2 sheet: insert AFP/LDH under pathology, shift hormones block, add HCG strata ---
$wsT2 = $wb.Worksheets.Item("Table 2")
$wsT2.Cells.Item(23,1).Value = 'pathology'
$wsT2.Cells.Item(23,2).Value = 'AFP, ng/mL'
$wsT2.Cells.Item(23,3).Value = '3 [IQR: 2 - 8.2]
range: 0.7 - 24000
complete: n = 436'
$wsT2.Cells.Item(24,1).Value = 'pathology'
$wsT2.Cells.Item(24,2).Value = 'LDH, U/L'
$wsT2.Cells.Item(24,3).Value = '200 [IQR: 180 - 260]
range: 3.1 - 2900
complete: n = 423'
$wsT2.Cells.Item(25,1).Value = 'pathology'
$wsT2.Cells.Item(25,2).Value = 'LDH strata'
$wsT2.Cells.Item(25,3).Value = '0 - 190 U/L: 39% (163)
> 190 U/L: 61% (260)
complete: n = 423'
$wsT2.Cells.Item(26,1).Value = 'pathology'
$wsT2.Cells.Item(26,2).Value = 'AFP strata'
$wsT2.Cells.Item(26,3).Value = '0 - 40 ng/mL: 85% (371)
> 40 ng/mL: 15% (65)
complete: n = 436'
$wsT2.Cells.Item(27,1).Value = 'hormones'
$wsT2.Cells.Item(27,2).Value = 'LH, mU/mL'
$wsT2.Cells.Item(27,3).Value = '3.2 [IQR: 1.4 - 5.6]
range: 0 - 47
complete: n = 370'
$wsT2.Cells.Item(28,1).Value = 'hormones'
$wsT2.Cells.Item(28,2).Value = 'FSH, mU/mL'
$wsT2.Cells.Item(28,3).Value = '5.2 [IQR: 1.7 - 9.7]
range: 0 - 100
complete: n = 368'
$wsT2.Cells.Item(29,1).Value = 'hormones'
$wsT2.Cells.Item(29,2).Value = 'PRL, µU/mL'
$wsT2.Cells.Item(29,3).Value = '120 [IQR: 14 - 200]
range: 2.4 - 1300
complete: n = 367'
$wsT2.Cells.Item(30,1).Value = 'hormones'
$wsT2.Cells.Item(30,2).Value = 'Total testosterone, ng/mL'
$wsT2.Cells.Item(30,3).Value = '4.4 [IQR: 3.3 - 6]
range: 1 - 17
complete: n = 374'
$wsT2.Cells.Item(31,1).Value = 'hormones'
$wsT2.Cells.Item(31,2).Value = 'Free testosterone, ng/mL'
$wsT2.Cells.Item(31,3).Value = '10 [IQR: 7.2 - 14]
range: 2.2 - 60
complete: n = 124'
$wsT2.Cells.Item(32,1).Value = 'hormones'
$wsT2.Cells.Item(32,2).Value = 'SHBG, nmol/L'
$wsT2.Cells.Item(32,3).Value = '32 [IQR: 22 - 44]
range: 6.3 - 160
complete: n = 39'
$wsT2.Cells.Item(33,1).Value = 'hormones'
$wsT2.Cells.Item(33,2).Value = 'HCG, IU/L'
$wsT2.Cells.Item(33,3).Value = '1 [IQR: 1 - 9]
range: 0.2 - 18000
complete: n = 436'
$wsT2.Cells.Item(34,1).Value = 'hormones'
$wsT2.Cells.Item(34,2).Value = 'E2, pg/mL'
$wsT2.Cells.Item(34,3).Value = '29 [IQR: 19 - 42]
range: 0.5 - 200
complete: n = 365'
$wsT2.Cells.Item(35,1).Value = 'hormones'
$wsT2.Cells.Item(35,2).Value = 'HCG strata'
$wsT2.Cells.Item(35,3).Value = '0 - 5 IU/L: 70% (306)
> 5 IU/L: 30% (130)
complete: n = 436'
$wsT2.Cells.Item(36,1).Value = 'treatment'
$wsT2.Cells.Item(36,2).Value = 'Type of surgery'
$wsT2.Cells.Item(36,3).Value = 'resection: 100% (438)
enucleation: 0.23% (1)
complete: n = 439'
$wsT2.Cells.Item(37,1).Value = 'treatment'
$wsT2.Cells.Item(37,2).Value = 'Chemotherapy'
$wsT2.Cells.Item(37,3).Value = '50% (219)
complete: n = 434'
$wsT2.Cells.Item(38,1).Value = 'treatment'
$wsT2.Cells.Item(38,2).Value = 'Radiation'
$wsT2.Cells.Item(38,3).Value = '3.9% (17)
complete: n = 433'
$wsT2.Cells.Item(39,1).Value = 'treatment'
$wsT2.Cells.Item(39,2).Value = 'Retroperitoneal lymphadenectomy'
$wsT2.Cells.Item(39,3).Value = '7.8% (34)
complete: n = 435'
$wsT2.Cells.Item(40,1).Value = 'treatment'
$wsT2.Cells.Item(40,2).Value = 'Testosterone replacement'
$wsT2.Cells.Item(40,3).Value = '8.4% (34)
complete: n = 406'
$wsT2.Cells.Item(41,1).Value = 'prognosis'
$wsT2.Cells.Item(41,2).Value = 'Follow-up, days'
$wsT2.Cells.Item(41,3).Value = '1300 [IQR: 540 - 2000]
range: 0 - 4700
complete: n = 439'
$wsT2.Cells.Item(42,1).Value = 'prognosis'
$wsT2.Cells.Item(42,2).Value = 'Relapse'
$wsT2.Cells.Item(42,3).Value = '7.3% (32)
complete: n = 439'

# --- Table 3 sheet: append LDH strata / AFP strata / HCG strata rows ---
$wsT3 = $wb.Worksheets.Item("Table 3")
$wsT3.Cells.Item(40,1).Value = 'LDH strata'
$wsT3.Cells.Item(40,2).Value = 423
$wsT3.Cells.Item(40,3).Value = 16
$wsT3.Cells.Item(40,4).Value = 3.644646924829157
$wsT3.Cells.Item(40,5).Value = 0.1470373942669988
$wsT3.Cells.Item(41,1).Value = 'AFP strata'
$wsT3.Cells.Item(41,2).Value = 436
$wsT3.Cells.Item(41,3).Value = 3
$wsT3.Cells.Item(41,4).Value = 0.683371298405467
$wsT3.Cells.Item(41,5).Value = 0.1106522587009888
$wsT3.Cells.Item(42,1).Value = 'HCG strata'
$wsT3.Cells.Item(42,2).Value = 436
$wsT3.Cells.Item(42,3).Value = 3
$wsT3.Cells.Item(42,4).Value = 0.683371298405467
$wsT3.Cells.Item(42,5).Value = 0.161569391982454

# --- Add new sheets Table 5 (relapse comparison) and Table 6 (LASSO models) ---
$sheetCount = $wb.Worksheets.Count
$wsT5 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($sheetCount))
$wsT5.Name = "Table 5"
$wsT6 = $wb.Worksheets.Add($null, $wsT5)
$wsT6.Name = "Table 6"

# Table 5 content
$wsT5.Cells.Item(1,1).Value = 'variable'
$wsT5.Cells.Item(1,2).Value = 'no_relapse'
$wsT5.Cells.Item(1,3).Value = 'relapse'
$wsT5.Cells.Item(1,4).Value = 'significance'
$wsT5.Cells.Item(1,5).Value = 'eff_size'
$wsT5.Cells.Item(2,1).Value = 'Age at surgery, years'
$wsT5.Cells.Item(2,2).Value = '35 [IQR: 28 - 43]
range: 17 - 86
complete: n = 407'
$wsT5.Cells.Item(2,3).Value = '35 [IQR: 27 - 42]
range: 17 - 57
complete: n = 32'
$wsT5.Cells.Item(2,4).Value = 'ns (p = 1)'
$wsT5.Cells.Item(2,5).Value = 'r = 0.0037'
$wsT5.Cells.Item(3,1).Value = 'Body mass index, kg/m2'
$wsT5.Cells.Item(3,2).Value = '25 [IQR: 23 - 28]
range: 17 - 41
complete: n = 306'
$wsT5.Cells.Item(3,3).Value = '25 [IQR: 23 - 28]
range: 21 - 32
complete: n = 22'
$wsT5.Cells.Item(3,4).Value = 'ns (p = 1)'
$wsT5.Cells.Item(3,5).Value = 'r = 0.006'
$wsT5.Cells.Item(4,1).Value = 'Body weight class'
$wsT5.Cells.Item(4,2).Value = 'normal: 47% (145)
overweight: 37% (114)
obesity: 15% (47)
complete: n = 306'
$wsT5.Cells.Item(4,3).Value = 'normal: 55% (12)
overweight: 27% (6)
obesity: 18% (4)
complete: n = 22'
$wsT5.Cells.Item(4,4).Value = 'ns (p = 0.98)'
$wsT5.Cells.Item(4,5).Value = 'V = 0.052'
$wsT5.Cells.Item(5,1).Value = 'Tumor stage'
$wsT5.Cells.Item(5,2).Value = 'I: 65% (264)
II: 30% (123)
III: 4.9% (20)
complete: n = 407'
$wsT5.Cells.Item(5,3).Value = 'I: 75% (24)
II: 22% (7)
III: 3.1% (1)
complete: n = 32'
$wsT5.Cells.Item(5,4).Value = 'ns (p = 0.95)'
$wsT5.Cells.Item(5,5).Value = 'V = 0.056'
$wsT5.Cells.Item(6,1).Value = 'Residual tumor'
$wsT5.Cells.Item(6,2).Value = 'R0: 98% (400)
R1: 1.7% (7)
complete: n = 407'
$wsT5.Cells.Item(6,3).Value = 'R0: 97% (31)
R1: 3.1% (1)
complete: n = 32'
$wsT5.Cells.Item(6,4).Value = 'ns (p = 1)'
$wsT5.Cells.Item(6,5).Value = 'V = 0.027'
$wsT5.Cells.Item(7,1).Value = 'Tumor size, cm'
$wsT5.Cells.Item(7,2).Value = '3 [IQR: 1.7 - 4.5]
range: 0.3 - 15
complete: n = 403'
$wsT5.Cells.Item(7,3).Value = '3.4 [IQR: 2.2 - 5]
range: 0.5 - 12
complete: n = 32'
$wsT5.Cells.Item(7,4).Value = 'ns (p = 0.88)'
$wsT5.Cells.Item(7,5).Value = 'r = 0.053'
$wsT5.Cells.Item(8,1).Value = 'Rete testis invaded'
$wsT5.Cells.Item(8,2).Value = '37% (148)
complete: n = 401'
$wsT5.Cells.Item(8,3).Value = '31% (10)
complete: n = 32'
$wsT5.Cells.Item(8,4).Value = 'ns (p = 0.98)'
$wsT5.Cells.Item(8,5).Value = 'V = 0.031'
$wsT5.Cells.Item(9,1).Value = 'Lymphovascular invasion'
$wsT5.Cells.Item(9,2).Value = '31% (126)
complete: n = 405'
$wsT5.Cells.Item(9,3).Value = '25% (8)
complete: n = 32'
$wsT5.Cells.Item(9,4).Value = 'ns (p = 0.98)'
$wsT5.Cells.Item(9,5).Value = 'V = 0.035'
$wsT5.Cells.Item(10,1).Value = 'Lugano class'
$wsT5.Cells.Item(10,2).Value = 'I: 80% (317)
II: 16% (63)
III: 4.5% (18)
complete: n = 398'
$wsT5.Cells.Item(10,3).Value = 'I: 94% (29)
II: 3.2% (1)
III: 3.2% (1)
complete: n = 31'
$wsT5.Cells.Item(10,4).Value = 'ns (p = 0.65)'
$wsT5.Cells.Item(10,5).Value = 'V = 0.095'
$wsT5.Cells.Item(11,1).Value = 'IGCCCG risk'
$wsT5.Cells.Item(11,2).Value = 'good: 95% (214)
intermediate: 4% (9)
poor: 1.3% (3)
complete: n = 226'
$wsT5.Cells.Item(11,3).Value = 'good: 100% (7)
intermediate: 0% (0)
poor: 0% (0)
complete: n = 7'
$wsT5.Cells.Item(11,4).Value = 'ns (p = 1)'
$wsT5.Cells.Item(11,5).Value = 'V = 0.041'
$wsT5.Cells.Item(12,1).Value = 'Histology'
$wsT5.Cells.Item(12,2).Value = 'seminoma: 64% (262)
mixed: 36% (145)
complete: n = 407'
$wsT5.Cells.Item(12,3).Value = 'seminoma: 55% (17)
mixed: 45% (14)
complete: n = 31'
$wsT5.Cells.Item(12,4).Value = 'ns (p = 0.95)'
$wsT5.Cells.Item(12,5).Value = 'V = 0.051'
$wsT5.Cells.Item(13,1).Value = 'Teratoma, %'
$wsT5.Cells.Item(13,2).Value = '0 [IQR: 0 - 0]
range: 0 - 100
complete: n = 404'
$wsT5.Cells.Item(13,3).Value = '0 [IQR: 0 - 7.5]
range: 0 - 100
complete: n = 31'
$wsT5.Cells.Item(13,4).Value = 'ns (p = 0.32)'
$wsT5.Cells.Item(13,5).Value = 'r = 0.096'
$wsT5.Cells.Item(14,1).Value = 'Embyonal cancer, %'
$wsT5.Cells.Item(14,2).Value = '0 [IQR: 0 - 10]
range: 0 - 100
complete: n = 404'
$wsT5.Cells.Item(14,3).Value = '0 [IQR: 0 - 3.8]
range: 0 - 95
complete: n = 30'
$wsT5.Cells.Item(14,4).Value = 'ns (p = 0.95)'
$wsT5.Cells.Item(14,5).Value = 'r = 0.03'
$wsT5.Cells.Item(15,1).Value = 'Chorion cancer, %'
$wsT5.Cells.Item(15,2).Value = '0 [IQR: 0 - 0]
range: 0 - 20
complete: n = 404'
$wsT5.Cells.Item(15,3).Value = '0 [IQR: 0 - 0]
range: 0 - 60
complete: n = 31'
$wsT5.Cells.Item(15,4).Value = 'ns (p = 0.65)'
$wsT5.Cells.Item(15,5).Value = 'r = 0.071'
$wsT5.Cells.Item(16,1).Value = 'Yolk sac cancer, %'
$wsT5.Cells.Item(16,2).Value = '0 [IQR: 0 - 0]
range: 0 - 100
complete: n = 404'
$wsT5.Cells.Item(16,3).Value = '0 [IQR: 0 - 0]
range: 0 - 100
complete: n = 30'
$wsT5.Cells.Item(16,4).Value = 'ns (p = 0.95)'
$wsT5.Cells.Item(16,5).Value = 'r = 0.041'
$wsT5.Cells.Item(17,1).Value = 'Seminoma, %'
$wsT5.Cells.Item(17,2).Value = '100 [IQR: 18 - 100]
range: 0 - 100
complete: n = 403'
$wsT5.Cells.Item(17,3).Value = '100 [IQR: 0 - 100]
range: 0 - 100
complete: n = 31'
$wsT5.Cells.Item(17,4).Value = 'ns (p = 0.79)'
$wsT5.Cells.Item(17,5).Value = 'r = 0.062'
$wsT5.Cells.Item(18,1).Value = 'Teratoma >= 75%'
$wsT5.Cells.Item(18,2).Value = '4.2% (17)
complete: n = 404'
$wsT5.Cells.Item(18,3).Value = '9.7% (3)
complete: n = 31'
$wsT5.Cells.Item(18,4).Value = 'ns (p = 0.95)'
$wsT5.Cells.Item(18,5).Value = 'V = 0.067'
$wsT5.Cells.Item(19,1).Value = 'Embryonic cancer >= 75%'
$wsT5.Cells.Item(19,2).Value = '12% (48)
complete: n = 404'
$wsT5.Cells.Item(19,3).Value = '3.3% (1)
complete: n = 30'
$wsT5.Cells.Item(19,4).Value = 'ns (p = 0.88)'
$wsT5.Cells.Item(19,5).Value = 'V = 0.069'
$wsT5.Cells.Item(20,1).Value = 'Yolk sac cancer >= 75%'
$wsT5.Cells.Item(20,2).Value = '1.2% (5)
complete: n = 404'
$wsT5.Cells.Item(20,3).Value = '3.3% (1)
complete: n = 30'
$wsT5.Cells.Item(20,4).Value = 'ns (p = 1)'
$wsT5.Cells.Item(20,5).Value = 'V = 0.046'
$wsT5.Cells.Item(21,1).Value = 'Seminoma >= 75%'
$wsT5.Cells.Item(21,2).Value = '68% (276)
complete: n = 403'
$wsT5.Cells.Item(21,3).Value = '61% (19)
complete: n = 31'
$wsT5.Cells.Item(21,4).Value = 'ns (p = 0.95)'
$wsT5.Cells.Item(21,5).Value = 'V = 0.04'
$wsT5.Cells.Item(22,1).Value = 'AFP, ng/mL'
$wsT5.Cells.Item(22,2).Value = '3 [IQR: 1.9 - 8.2]
range: 0.7 - 24000
complete: n = 404'
$wsT5.Cells.Item(22,3).Value = '3 [IQR: 2.3 - 7.6]
range: 1.2 - 9400
complete: n = 32'
$wsT5.Cells.Item(22,4).Value = 'ns (p = 0.95)'
$wsT5.Cells.Item(22,5).Value = 'r = 0.038'
$wsT5.Cells.Item(23,1).Value = 'LDH, U/L'
$wsT5.Cells.Item(23,2).Value = '200 [IQR: 180 - 260]
range: 3.1 - 2900
complete: n = 395'
$wsT5.Cells.Item(23,3).Value = '200 [IQR: 170 - 260]
range: 120 - 1000
complete: n = 28'
$wsT5.Cells.Item(23,4).Value = 'ns (p = 0.98)'
$wsT5.Cells.Item(23,5).Value = 'r = 0.02'
$wsT5.Cells.Item(24,1).Value = 'LDH strata'
$wsT5.Cells.Item(24,2).Value = '0 - 190 U/L: 38% (151)
> 190 U/L: 62% (244)
complete: n = 395'
$wsT5.Cells.Item(24,3).Value = '0 - 190 U/L: 43% (12)
> 190 U/L: 57% (16)
complete: n = 28'
$wsT5.Cells.Item(24,4).Value = 'ns (p = 1)'
$wsT5.Cells.Item(24,5).Value = 'V = 0.024'
$wsT5.Cells.Item(25,1).Value = 'AFP strata'
$wsT5.Cells.Item(25,2).Value = '0 - 40 ng/mL: 85% (345)
> 40 ng/mL: 15% (59)
complete: n = 404'
$wsT5.Cells.Item(25,3).Value = '0 - 40 ng/mL: 81% (26)
> 40 ng/mL: 19% (6)
complete: n = 32'
$wsT5.Cells.Item(25,4).Value = 'ns (p = 0.98)'
$wsT5.Cells.Item(25,5).Value = 'V = 0.03'
$wsT5.Cells.Item(26,1).Value = 'LH, mU/mL'
$wsT5.Cells.Item(26,2).Value = '3.2 [IQR: 1.4 - 5.6]
range: 0 - 47
complete: n = 346'
$wsT5.Cells.Item(26,3).Value = '3.3 [IQR: 1.4 - 6.7]
range: 0 - 25
complete: n = 24'
$wsT5.Cells.Item(26,4).Value = 'ns (p = 1)'
$wsT5.Cells.Item(26,5).Value = 'r = 0.00015'
$wsT5.Cells.Item(27,1).Value = 'FSH, mU/mL'
$wsT5.Cells.Item(27,2).Value = '5.2 [IQR: 1.7 - 9.6]
range: 0 - 100
complete: n = 344'
$wsT5.Cells.Item(27,3).Value = '4.8 [IQR: 2.2 - 11]
range: 0 - 38
complete: n = 24'
$wsT5.Cells.Item(27,4).Value = 'ns (p = 1)'
$wsT5.Cells.Item(27,5).Value = 'r = 0.006'
$wsT5.Cells.Item(28,1).Value = 'PRL, µU/mL'
$wsT5.Cells.Item(28,2).Value = '120 [IQR: 13 - 200]
range: 2.4 - 1300
complete: n = 343'
$wsT5.Cells.Item(28,3).Value = '170 [IQR: 130 - 210]
range: 86 - 620
complete: n = 24'
$wsT5.Cells.Item(28,4).Value = 'p = 0.039'
$wsT5.Cells.Item(28,5).Value = 'r = 0.15'
$wsT5.Cells.Item(29,1).Value = 'Total testosterone, ng/mL'
$wsT5.Cells.Item(29,2).Value = '4.4 [IQR: 3.3 - 6]
range: 1 - 17
complete: n = 349'
$wsT5.Cells.Item(29,3).Value = '4.2 [IQR: 3.6 - 5.7]
range: 1.8 - 10
complete: n = 25'
$wsT5.Cells.Item(29,4).Value = 'ns (p = 1)'
$wsT5.Cells.Item(29,5).Value = 'r = 0.015'
$wsT5.Cells.Item(30,1).Value = 'HCG, IU/L'
$wsT5.Cells.Item(30,2).Value = '1 [IQR: 1 - 7.1]
range: 0.2 - 18000
complete: n = 404'
$wsT5.Cells.Item(30,3).Value = '2 [IQR: 1 - 50]
range: 1 - 5600
complete: n = 32'
$wsT5.Cells.Item(30,4).Value = 'ns (p = 0.19)'
$wsT5.Cells.Item(30,5).Value = 'r = 0.11'
$wsT5.Cells.Item(31,1).Value = 'E2, pg/mL'
$wsT5.Cells.Item(31,2).Value = '29 [IQR: 18 - 42]
range: 0.5 - 200
complete: n = 341'
$wsT5.Cells.Item(31,3).Value = '32 [IQR: 21 - 41]
range: 0.5 - 130
complete: n = 24'
$wsT5.Cells.Item(31,4).Value = 'ns (p = 0.98)'
$wsT5.Cells.Item(31,5).Value = 'r = 0.022'
$wsT5.Cells.Item(32,1).Value = 'HCG strata'
$wsT5.Cells.Item(32,2).Value = '0 - 5 IU/L: 71% (288)
> 5 IU/L: 29% (116)
complete: n = 404'
$wsT5.Cells.Item(32,3).Value = '0 - 5 IU/L: 56% (18)
> 5 IU/L: 44% (14)
complete: n = 32'
$wsT5.Cells.Item(32,4).Value = 'ns (p = 0.65)'
$wsT5.Cells.Item(32,5).Value = 'V = 0.086'
$wsT5.Cells.Item(33,1).Value = 'Type of surgery'
$wsT5.Cells.Item(33,2).Value = 'resection: 100% (406)
enucleation: 0.25% (1)
complete: n = 407'
$wsT5.Cells.Item(33,3).Value = 'resection: 100% (32)
enucleation: 0% (0)
complete: n = 32'
$wsT5.Cells.Item(33,4).Value = 'ns (p = 1)'
$wsT5.Cells.Item(33,5).Value = 'V = 0.013'
$wsT5.Cells.Item(34,1).Value = 'Chemotherapy'
$wsT5.Cells.Item(34,2).Value = '53% (212)
complete: n = 402'
$wsT5.Cells.Item(34,3).Value = '22% (7)
complete: n = 32'
$wsT5.Cells.Item(34,4).Value = 'p = 0.027'
$wsT5.Cells.Item(34,5).Value = 'V = 0.16'
$wsT5.Cells.Item(35,1).Value = 'Radiation'
$wsT5.Cells.Item(35,2).Value = '4.2% (17)
complete: n = 401'
$wsT5.Cells.Item(35,3).Value = '0% (0)
complete: n = 32'
$wsT5.Cells.Item(35,4).Value = 'ns (p = 0.95)'
$wsT5.Cells.Item(35,5).Value = 'V = 0.057'
$wsT5.Cells.Item(36,1).Value = 'Retroperitoneal lymphadenectomy'
$wsT5.Cells.Item(36,2).Value = '8.2% (33)
complete: n = 403'
$wsT5.Cells.Item(36,3).Value = '3.1% (1)
complete: n = 32'
$wsT5.Cells.Item(36,4).Value = 'ns (p = 0.95)'
$wsT5.Cells.Item(36,5).Value = 'V = 0.049'
$wsT5.Cells.Item(37,1).Value = 'Testosterone replacement'
$wsT5.Cells.Item(37,2).Value = '7% (26)
complete: n = 374'
$wsT5.Cells.Item(37,3).Value = '25% (8)
complete: n = 32'
$wsT5.Cells.Item(37,4).Value = 'p = 0.027'
$wsT5.Cells.Item(37,5).Value = 'V = 0.18'
$wsT5.Rows.Item(1).Font.Bold = $true
$wsT5.Rows.Item(1).HorizontalAlignment = -4108

# Table 6 content
$wsT6.Cells.Item(1,1).Value = 'LASSO model type'
$wsT6.Cells.Item(1,2).Value = 'Dataset'
$wsT6.Cells.Item(1,3).Value = 'Concordance index'
$wsT6.Cells.Item(1,4).Value = 'Nagelkirke R²'
$wsT6.Cells.Item(1,5).Value = 'Integrated Brier Score'
$wsT6.Cells.Item(2,1).Value = 'fist-/second-term PRL'
$wsT6.Cells.Item(2,2).Value = 'data'
$wsT6.Cells.Item(2,3).Value = 0.85
$wsT6.Cells.Item(2,4).Value = 0.26
$wsT6.Cells.Item(2,5).Value = 0.048
$wsT6.Cells.Item(3,1).Value = 'fist-/second-term PRL'
$wsT6.Cells.Item(3,2).Value = '10-fold cross-validation'
$wsT6.Cells.Item(3,3).Value = 0.84
$wsT6.Cells.Item(3,4).Value = 0.25
$wsT6.Cells.Item(3,5).Value = 0.049
$wsT6.Cells.Item(4,1).Value = 'PRL strata'
$wsT6.Cells.Item(4,2).Value = 'data'
$wsT6.Cells.Item(4,3).Value = 0.86
$wsT6.Cells.Item(4,4).Value = 0.29
$wsT6.Cells.Item(4,5).Value = 0.049
$wsT6.Cells.Item(5,1).Value = 'PRL strata'
$wsT6.Cells.Item(5,2).Value = '10-fold cross-validation'
$wsT6.Cells.Item(5,3).Value = 0.85
$wsT6.Cells.Item(5,4).Value = 0.27
$wsT6.Cells.Item(5,5).Value = 0.05
$wsT6.Rows.Item(1).Font.Bold = $true
$wsT6.Rows.Item(1).HorizontalAlignment = -4108

# --- Restore original active sheet (Cover) ---
$wsCover.Activate()
$wsCover.Range("A1").Select() | Out-Null
